# Add the team's season record (Wins / Losses / Ties) as three new
# columns (AD, AE, AF) to the right of the existing player-stats table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, styled like the rest of the header row (row 1).
$ws.Range("AD1").Value2 = "Wins"
$ws.Range("AE1").Value2 = "Losses"
$ws.Range("AF1").Value2 = "Ties"

# Reuse the same header formatting (bold font + border + centered/top
# alignment) already applied to A1:AC1, so the new header cells share
# the exact same cell style as their neighbours.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Season record for every player row (2-50): 75 wins, 86 losses, 0 ties.
$wins = 75
$losses = 86
$ties = 0

for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value2 = $wins
    $ws.Cells.Item($r, 31).Value2 = $losses
    $ws.Cells.Item($r, 32).Value2 = $ties
}
